# The library's default `fieldMatchType` changed to `labelTypeBrackets`, so
# the generated header row now combines each column's label with its type
# in brackets instead of the bare label:
#   "Order Number" -> "Order Number[OrderNo]"
#   "User ID"      -> "User ID[buyer]"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Order Number[OrderNo]"
$ws.Range("B1").Value = "User ID[buyer]"

# Match the saved selection in the updated sample file (header row selected
# instead of the whole column C).
$ws.Range("A1:B1").Select() | Out-Null
